# Updated cryptos list on Sat Jan 13 16:18:51 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '43.075.79'
$ws.Range("E2").Value = '  -3.06%  '
$ws.Range("D3").Value = '2.562.76'
$ws.Range("E3").Value = '  -3.80%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.83'
$ws.Range("E5").Value = '  -1.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.65'
$ws.Range("E6").Value = '  -4.41%  '
$ws.Range("E7").Value = '  -2.52%  '
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.549'
$ws.Range("E9").Value = '  -3.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.30'
$ws.Range("E10").Value = '  -3.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0812'
$ws.Range("E11").Value = '  -1.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.67'
$ws.Range("E12").Value = '  -3.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.114'
$ws.Range("E13").Value = '  +6.99%  '
$ws.Range("D14").Value = '2.952.77'
$ws.Range("E14").Value = '  -3.97%  '
$ws.Range("D15").Value = '2.624.40'
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.882'
$ws.Range("E16").Value = '  -3.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.21'
$ws.Range("E17").Value = '  -4.34%  '
$ws.Range("D18").Value = '43.050.52'
$ws.Range("E18").Value = '  -3.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.75'
$ws.Range("E19").Value = '  +1.68%  '
$ws.Range("D20").Value = '0.0₃0984'
$ws.Range("E20").Value = '  -1.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.55'
$ws.Range("E21").Value = '  -3.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.85'
$ws.Range("E22").Value = '  -3.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '253.94'
$ws.Range("E23").Value = '  -7.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.95'
$ws.Range("E24").Value = '  -1.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.12'
$ws.Range("E25").Value = '  -7.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '29.01'
$ws.Range("E26").Value = '  -4.34%  '
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.24'
$ws.Range("E28").Value = '  -1.17%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.21'
$ws.Range("E29").Value = '  +0.85%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.04'
$ws.Range("E30").Value = '  -1.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.14'
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '153.18'
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.76'
$ws.Range("E33").Value = '  -1.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.38'
$ws.Range("E34").Value = '  -8.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.12'
$ws.Range("E35").Value = '  -8.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0796'
$ws.Range("E36").Value = '  -3.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.114'
$ws.Range("E37").Value = '  -3.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.60'
$ws.Range("E38").Value = '  +11.62%  '
$ws.Range("E39").Value = '  -2.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.67'
$ws.Range("E40").Value = '  -9.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.14'
$ws.Range("E41").Value = '  +34.80%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0312'
$ws.Range("E42").Value = '  -2.38%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.42'
$ws.Range("E43").Value = '  -3.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.86'
$ws.Range("E44").Value = '  -0.94%  '
$ws.Range("D45").Value = '2.101.97'
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.997'
$ws.Range("E46").Value = '  -0.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.22'
$ws.Range("E47").Value = '  -1.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.31'
$ws.Range("E48").Value = '  -6.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.73'
$ws.Range("E49").Value = '  -3.47%  '
$ws.Range("D50").Value = '2.808.40'
$ws.Range("E50").Value = '  -3.88%  '
$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.85'
$ws.Range("E51").Value = '  +5.00%  '
